$d = $word.ActiveDocument

# The document currently starts with:
#   Paragraph 1 (style Heading1): "Untitled Review of 'A Christian Looks at the Jewish Question'"
#   Paragraph 2 (no style, bold run): "By Dorothy Day"
#
# We need to turn this into a pandoc-style title block:
#   Paragraph 1 (style Title): "Untitled Review of 'A Christian Looks at the Jewish Question'"
#     split into one run per word/space/apostrophe token
#   Paragraph 2 (style Authors): "Dorothy Day"
#     split into "Dorothy" / " " / "Day" runs, no bold

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function New-RunXml([string]$text) {
    $escaped = $text.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    return "<w:r $wNs><w:t xml:space=`"preserve`">$escaped</w:t></w:r>"
}

$titleTokens = @("Untitled", " ", "Review", " ", "of", " ", "'", "A", " ", "Christian", " ", "Looks", " ", "at", " ", "the", " ", "Jewish", " ", "Question", "'")
$authorTokens = @("Dorothy", " ", "Day")

$titleRuns = ($titleTokens | ForEach-Object { New-RunXml $_ }) -join ""
$authorRuns = ($authorTokens | ForEach-Object { New-RunXml $_ }) -join ""

$titleParaXml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"Title`"/></w:pPr>$titleRuns</w:p>"
$authorParaXml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"Authors`"/></w:pPr>$authorRuns</w:p>"

$p1 = $d.Paragraphs(1)
$p2 = $d.Paragraphs(2)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
[void]$rng.InsertXML($titleParaXml + $authorParaXml)

Write-Host "Title/author block rewritten."
